# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.987.40"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "2.552.06"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.99"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.28"
$ws.Range("E6").Value = "  +3.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.578"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.80"
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0824"
$ws.Range("E11").Value = "  +2.01%  "
$ws.Range("E12").Value = "  +5.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.65"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").Value = "2.946.02"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "2.510.82"
$ws.Range("E15").Value = "  -3.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.98"
$ws.Range("E16").Value = "  +6.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.875"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "43.023.95"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.79"
$ws.Range("E19").Value = "  +4.32%  "
$ws.Range("D20").Value = "0.0₃0997"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.62"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.03"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.87"
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("E25").Value = "  -2.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.20"
$ws.Range("E26").Value = "  -4.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.21"
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.16"
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.15"
$ws.Range("E31").Value = "  +2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.95"
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.49"
$ws.Range("E35").Value = "  +14.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0806"
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.32"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.117"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.68"
$ws.Range("E39").Value = "  +10.34%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.45"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.09"
$ws.Range("E42").Value = "  +31.63%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("D45").Value = "2.086.42"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.85"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.00"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").Value = "2.803.29"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.03"
$ws.Range("E50").Value = "  +8.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.62"
$ws.Range("E51").Value = "  -1.67%  "
